# Fixed issue with modifying uploaded images
#
# 1. Rename the merge-field template tags from "images" to "images_buffer"
#    (both the opening "{#images}" and closing "{/images}" tags).
# 2. Move the "_GoBack" bookmark from right after the opening tag to right
#    before the closing "}" of the closing tag.
# 3. Refresh the cached header date field result.

$d = $word.ActiveDocument

# --- 1a. Opening tag: {#images} -> {#images_buffer} ---------------------
$d.Content.Find.Execute("{#images}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{#images_buffer}", 2) | Out-Null

# --- 1b. Closing tag: {/images} -> {/images_buffer} ----------------------
$d.Content.Find.Execute("{/images}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{/images_buffer}", 2) | Out-Null

# --- 2. Move the _GoBack bookmark to sit just before the final "}" -------
$rng = $d.Content
$rng.Find.Execute("{/images_buffer}", $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$bmPos = $rng.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 3. Refresh cached header date field result ---------------------------
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("10/28/18", $false, $false, $false, $false, $false, `
                                     $true, 1, $false, "11/4/18", 2) | Out-Null
        }
    }
}
